$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$letters = @("a", "b", "c", "d", "e", "f", "g", "h", "i", "j")

for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = $i + 1
    $letter = $letters[$i]
    $ws.Cells.Item($row, 2).Value = $letter
    $ws.Cells.Item($row, 3).Value = $letter
    $ws.Cells.Item($row, 4).Value = $letter
    $ws.Cells.Item($row, 5).Value = $letter
}

$ws.Range("B12").Select()
